# Importing of Customers and Vendors working
#
# The "terms" column (P) becomes a numeric sequence number (1, 2, ...)
# and the old "yes" flag column (Q) becomes an explicit YES/NO per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (first vendor record)
$ws.Range("P1").Value = 1
$ws.Range("Q1").Value = "YES"

# Row 2 (second vendor record)
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = "NO"

# Scroll the view over to column P and select P3, matching the saved
# window state after reviewing the newly imported columns.
$ws.Application.ActiveWindow.ScrollColumn = 16
$ws.Range("P3").Select()
